$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (they reference old content that is being replaced)
$ws.Hyperlinks.Delete()

# Row 1: plain text header (no hyperlink)
$ws.Range("A1").Value = "Checked and posting"

# Rows 2-10: text + hyperlink pairs
$rows = @(
    @{ Cell = "A2";  Text = "Diabetic Food and Recipes | Facebook"; Url = "https://www.facebook.com/groups/463604924157606" },
    @{ Cell = "A3";  Text = "Diabetic Support Group | Facebook"; Url = "https://www.facebook.com/groups/144413888950964" },
    @{ Cell = "A4";  Text = "Healthy food planners for diabetic people. Planeamiento de alimentacion | Facebook"; Url = "https://www.facebook.com/groups/197360433413" },
    @{ Cell = "A5";  Text = "Diabetic Recipes New Ideas 2024 | Facebook"; Url = "https://www.facebook.com/groups/1006148591253340" },
    @{ Cell = "A6";  Text = "DAILY DIABETIC DIET RECIPES | Facebook"; Url = "https://www.facebook.com/groups/652938943691405" },
    @{ Cell = "A7";  Text = "Easy Recipes Group | Facebook"; Url = "https://www.facebook.com/groups/easyrecipes105/" },
    @{ Cell = "A8";  Text = "Diabetes Type 2 Signs And Symptoms, Foods, Recipes, Meal, Support Group USA | Facebook"; Url = "https://www.facebook.com/groups/diabetestype2signsandsymptomsfoodsrecipesmeal/" },
    @{ Cell = "A9";  Text = "Recipes Taught | Facebook"; Url = "https://www.facebook.com/groups/recipestaught/" },
    @{ Cell = "A10"; Text = "Easy Diabetic Recipes - cooking | foods | Facebook"; Url = "https://www.facebook.com/groups/foods.cooking.recipes/" }
)

foreach ($row in $rows) {
    $cell = $ws.Range($row.Cell)
    $ws.Hyperlinks.Add($cell, $row.Url, [Type]::Missing, [Type]::Missing, $row.Url)
    $cell.Value = $row.Text
    $cell.Style = "Hyperlink"
}

$ws.Range("A11").Select()
